$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (characters) so the saved OOXML column width becomes 15
# (ColumnWidth uses Excel's internal character-width metric, which does not
# map 1:1 onto the OOXML "width" attribute; 14.17 round-trips to width="15").
$ws.Columns.Item(1).ColumnWidth = 14.17

# Update the Priority (C) and Urgency (D) values for rows 2-11
$ws.Range("C2").Value = 402.04
$ws.Range("D2").Value = 350.3

$ws.Range("C3").Value = 128.15
$ws.Range("D3").Value = 114.1

$ws.Range("C4").Value = 53.04
$ws.Range("D4").Value = 55.28

$ws.Range("C5").Value = 75.11
$ws.Range("D5").Value = 58.82

$ws.Range("C6").Value = 114.85
$ws.Range("D6").Value = 113.28

$ws.Range("C7").Value = 57.11
$ws.Range("D7").Value = 56.05

$ws.Range("C8").Value = 57.75
$ws.Range("D8").Value = 57.22

$ws.Range("C9").Value = 159.04
$ws.Range("D9").Value = 122.93

$ws.Range("C10").Value = 82.6
$ws.Range("D10").Value = 65.98

$ws.Range("C11").Value = 76.44
$ws.Range("D11").Value = 56.95
